$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '34.035.74'
$ws.Range("E2").Value = '  -1.58%  '

$ws.Range("D3").Value = '1.793.49'
$ws.Range("E3").Value = '  -1.56%  '

$ws.Range("E4").Value = '  +0.72%  '

Set-TextValue $ws.Range("D5") '227.73'
$ws.Range("E5").Value = '  -2.84%  '

Set-TextValue $ws.Range("D6") '0.555'
$ws.Range("E6").Value = '  +0.55%  '

$ws.Range("E7").Value = '  +0.68%  '

Set-TextValue $ws.Range("D8") '31.11'
$ws.Range("E8").Value = '  -2.82%  '

Set-TextValue $ws.Range("D9") '46.17'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("E11").Value = '  -3.29%  '

$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("D13").Value = '2.051.46'
$ws.Range("E13").Value = '  -1.09%  '

Set-TextValue $ws.Range("D14") '11.26'
$ws.Range("E14").Value = '  +8.69%  '

$ws.Range("D15").Value = '1.813.32'
$ws.Range("E15").Value = '  -0.73%  '

Set-TextValue $ws.Range("D16") '0.635'
$ws.Range("E16").Value = '  -1.97%  '

$ws.Range("D17").Value = '34.071.11'
$ws.Range("E17").Value = '  -1.25%  '

$ws.Range("E18").Value = '  -3.38%  '

Set-TextValue $ws.Range("D19") '69.65'
$ws.Range("E19").Value = '  -2.74%  '

Set-TextValue $ws.Range("D20") '253.33'
$ws.Range("E20").Value = '  -4.52%  '

$ws.Range("D21").Value = '0.0₃0746'
$ws.Range("E21").Value = '  -2.61%  '

$ws.Range("E22").Value = '  +0.50%  '

Set-TextValue $ws.Range("D23") '10.45'
$ws.Range("E23").Value = '  -0.90%  '

$ws.Range("E24").Value = '  -3.42%  '

$ws.Range("E25").Value = '  -1.65%  '

Set-TextValue $ws.Range("D26") '157.88'
$ws.Range("E26").Value = '  -2.76%  '

Set-TextValue $ws.Range("D27") '16.62'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D28") '7.03'
$ws.Range("E28").Value = '  -2.14%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D29") '0.114'
$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("E30").Value = '  +0.74%  '

Set-TextValue $ws.Range("D31") '3.90'
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("E33").Value = '  -0.58%  '

$ws.Range("E34").Value = '  +0.82%  '

$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").Value = '1.485.47'
$ws.Range("E36").Value = '  -6.26%  '

Set-TextValue $ws.Range("D37") '1.07'
$ws.Range("E37").Value = '  -0.07%  '

Set-TextValue $ws.Range("D38") '0.636'
$ws.Range("E38").Value = '  +0.67%  '

Set-TextValue $ws.Range("D39") '0.0187'
$ws.Range("E39").Value = '  -1.07%  '

Set-TextValue $ws.Range("D40") '83.87'
$ws.Range("E40").Value = '  -6.60%  '

Set-TextValue $ws.Range("D41") '2.83'
$ws.Range("E41").Value = '  -1.22%  '

$ws.Range("E42").Value = '  -0.29%  '

Set-TextValue $ws.Range("D43") '0.906'
$ws.Range("E43").Value = '  -2.69%  '

Set-TextValue $ws.Range("D44") '2.07'
$ws.Range("E44").Value = '  -4.21%  '

$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("E46").Value = '  +1.98%  '

$ws.Range("D47").Value = '1.949.95'
$ws.Range("E47").Value = '  -0.22%  '

Set-TextValue $ws.Range("D48") '5.72'
$ws.Range("E48").Value = '  -1.49%  '

$ws.Range("E49").Value = '  +0.38%  '

Set-TextValue $ws.Range("D50") '11.80'
$ws.Range("E50").Value = '  +2.11%  '

Set-TextValue $ws.Range("D51") '51.46'
$ws.Range("E51").Value = '  -5.16%  '
